# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# -> the recomputed "K" column values (column G) for each trade row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 0
    6  = 1
    7  = 0
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 2
    13 = 0
    14 = 2
    15 = 1
    16 = 0
    17 = 1
    18 = 1
    19 = 2
    20 = 2
    21 = 2
    22 = 2
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    27 = 0
    28 = 1
    29 = 0
    31 = 3
    33 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
